$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("G900153336_calls")

# Copy formatting of an existing fully-styled data row (row 2) onto the new row 13
$ws.Range("A2:E2").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row 13 values in the order the strings were first introduced
$ws.Range("C13").Value = "Rejected"
$ws.Range("F13").Value = "Verify that rejected specimen's status is updated and has no result created"
$ws.Range("A13").Value = "SpecimenAntibodyResults10"

# Column F widened to fit the longer comment text (closest reproducible width
# to the authored 65.6640625 given this engine's internal rounding)
$ws.Columns.Item(6).ColumnWidth = 64.75

# Page orientation explicitly set to portrait
$ws.PageSetup.Orientation = 1

# Update the active selection to E13
[void]$ws.Range("E13").Select()

Write-Host "Applied CVDLS-296 rejected specimen test row"
